# V0.5.0 Updated blogs and home. Some additional information about the application context
$wb = $excel.ActiveWorkbook

# --- faq sheet: fix grammar in the FAQ question (row 10, column A) ---
$faq = $wb.Worksheets.Item("faq")
$faq.Range("A10").Value = "Why a bean call is better than a remote method call?"

# --- links sheet: append a new link/reference row ---
$links = $wb.Worksheets.Item("links")
$links.Range("A13").Value = "https://www.sueddeutsche.de/wissen/kuenstliche-intelligenz-software-computer-1.5036926?utm_source=pocket-newtab-global-de-DE"

# --- view state: move the active selection on faq down a few rows, and
#     make "links" (with the newly added row) the active sheet/selection ---
$faq.Activate()
$faq.Range("B10").Select()

$links.Activate()
$links.Range("A13").Select()
